# Update the "Trace_Report_SBHP_initial" trace-report sheet with the new
# (re-run) search results: 3 events instead of 2, new car numbers/weights,
# and an extra data row. Also drop the AutoFilter / FilterDatabase defined
# name that no longer applies to the refreshed data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: search-summary banner text -------------------------------------
$ws.Range("A1").Value = "Description unknown, completed 06/15/2023 05:55:31 EDT, by WPJTOWN1.The search returned: 3 events."

# --- Row 3: car CGAX10167 (was BNGX30307) -----------------------------------
$ws.Range("A3").Value = "CGAX"
$ws.Range("B3").Value = 10167
$ws.Range("C3").Value = "GREELEY"
$ws.Range("D3").Value = "CO"
$ws.Range("E3").Value = 6
$ws.Range("F3").Value = 9
$ws.Range("G3").Value = 1217
$ws.Range("H3").Value = "Junction Delivery"
$ws.Range("I3").Value = "GWR"
$ws.Range("J3").Value = "JOHNSTOWN"
$ws.Range("K3").Value = "CO"
$ws.Range("L3").Value = 273000
$ws.Range("M3").Value = 64200
$ws.Range("N3").Value = 208800
$ws.Range("O3").Value = "CGAX10167"

# --- Row 4: car BNGX30727 (was UP74772) -------------------------------------
$ws.Range("A4").Value = "BNGX"
$ws.Range("B4").Value = 30727
$ws.Range("C4").Value = "KANSAS CITY"
$ws.Range("D4").Value = "KS"
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 14
$ws.Range("G4").Value = 853
$ws.Range("H4").Value = "Arrive In-Transit"
$ws.Range("I4").Value = "HLINKC"
$ws.Range("J4").Value = "LOVELAND"
$ws.Range("K4").Value = "CO"
$ws.Range("L4").Value = 282200
$ws.Range("M4").Value = 64400
$ws.Range("N4").Value = 217800
$ws.Range("O4").Value = "BNGX30727"

# --- Row 5 (new): car CGEX1941 ----------------------------------------------
$ws.Range("A5").Value = "CGEX"
$ws.Range("B5").Value = 1941
$ws.Range("C5").Value = "SHORTLINE YARD"
$ws.Range("D5").Value = "IA"
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 14
$ws.Range("G5").Value = 1913
$ws.Range("H5").Value = "Arrive In-Transit"
$ws.Range("I5").Value = "LTJ40N"
$ws.Range("J5").Value = "JOHNSTOWN"
$ws.Range("K5").Value = "CO"
$ws.Range("L5").Value = 198750
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 198750
$ws.Range("O5").Value = "CGEX1941"

# --- Drop the AutoFilter (and its _FilterDatabase defined name) ------------
$ws.AutoFilterMode = $false
foreach ($n in $wb.Names) {
    $n.Delete()
}

# --- Selection now spans the extra row --------------------------------------
$ws.Range("O3:O5").Select()
